$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.771.24"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.704.24"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D5").Value = "317.16"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.3942"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "1.525"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "53.67"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "0.08904"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "7.509"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").Value = "23.74"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "8.156"
$ws.Range("E15").Value = "  +6.60%  "
$ws.Range("D16").Value = "0.00001329"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "1.712.96"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "100.07"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").Value = "0.07054"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "19.76"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "7.100"
$ws.Range("E21").Value = "  +2.99%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +2.14%  "
$ws.Range("D24").Value = "24.761.20"
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("D25").Value = "3.242"
$ws.Range("E25").Value = "  +4.98%  "
$ws.Range("D26").Value = "2.370"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").Value = "22.80"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").Value = "162.53"
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "8.834"
$ws.Range("E29").Value = "  +15.64%  "
$ws.Range("D30").Value = "136.54"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "5.180"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "7.978"
$ws.Range("E32").Value = "  +5.84%  "
$ws.Range("D33").Value = "0.08943"
$ws.Range("E33").Value = "  +3.80%  "
$ws.Range("D34").Value = "1.086"
$ws.Range("E34").Value = "  -3.29%  "
$ws.Range("D35").Value = "1.983"
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "0.2766"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").Value = "0.02793"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "0.09199"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").Value = "1.464"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").Value = "0.7733"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "15.88"
$ws.Range("E43").Value = "  +2.13%  "
$ws.Range("D45").Value = "2.577"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("D46").Value = "4.214"
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").Value = "1.342"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "141.04"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "90.92"
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").Value = "0.07986"
$ws.Range("E51").Value = "  -0.90%  "
